$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Hbegf", "Erbb2", "ECs", 3, 1, 7.379281, 22.137843, 0.2744121884499962, 0.2744121884499961, 2, 0.6666666666666666, 1.720171333333333, 5.160514, 0.1961456356393658, 0.1961456356393658, 12.69362763681133, 114.242648731302, 0.05382475313071392, 0.05382475313071391),
    @("ECs", "Hbegf", "Erbb2", "FAPs", 3, 1, 7.379281, 22.137843, 0.2744121884499962, 0.2744121884499961, 3, 1, 4.076388666666666, 12.229166, 0.464817562438416, 0.464817562438416, 30.08081743654866, 270.727356928938, 0.1275516045387185, 0.1275516045387184),
    @("ECs", "Hbegf", "Erbb2", "sCs", 3, 1, 7.379281, 22.137843, 0.2744121884499962, 0.2744121884499961, 3, 1, 2.973308, 8.919924, 0.3390368019222182, 0.3390368019222182, 21.940875231548, 197.467877083932, 0.09303583078056377, 0.09303583078056374),
    @("FAPs", "Hbegf", "Erbb2", "ECs", 3, 1, 13.29805733333333, 39.894172, 0.4945128143207339, 0.4945128143207338, 2, 0.6666666666666666, 1.720171333333333, 5.160514, 0.1961456356393658, 0.1961456356393658, 22.87493701382311, 205.874433124408, 0.09699653029675201, 0.096996530296752),
    @("FAPs", "Hbegf", "Erbb2", "FAPs", 3, 1, 13.29805733333333, 39.894172, 0.4945128143207339, 0.4945128143207338, 3, 1, 4.076388666666666, 12.229166, 0.464817562438416, 0.464817562438416, 54.20805020228355, 487.8724518205519, 0.2298582409471246, 0.2298582409471245),
    @("FAPs", "Hbegf", "Erbb2", "sCs", 3, 1, 13.29805733333333, 39.894172, 0.4945128143207339, 0.4945128143207338, 3, 1, 2.973308, 8.919924, 0.3390368019222182, 0.3390368019222182, 39.53922025365866, 355.852982282928, 0.1676580430768574, 0.1676580430768573),
    @("M2", "Hbegf", "Erbb2", "ECs", 3, 1, 1.673314666666667, 5.019944000000001, 0.06222529534320158, 0.06222529534320156, 2, 0.6666666666666666, 1.720171333333333, 5.160514, 0.1961456356393658, 0.1961456356393658, 2.878387921246222, 25.905491291216, 0.01220522010793954, 0.01220522010793954),
    @("M2", "Hbegf", "Erbb2", "FAPs", 3, 1, 1.673314666666667, 5.019944000000001, 0.06222529534320158, 0.06222529534320156, 3, 1, 4.076388666666666, 12.229166, 0.464817562438416, 0.464817562438416, 6.821080942967112, 61.389728486704, 0.02892341010343748, 0.02892341010343747),
    @("M2", "Hbegf", "Erbb2", "sCs", 3, 1, 1.673314666666667, 5.019944000000001, 0.06222529534320158, 0.06222529534320156, 3, 1, 2.973308, 8.919924, 0.3390368019222182, 0.3390368019222182, 4.975279884917334, 44.777518964256, 0.02109666513182456, 0.02109666513182456),
    @("sCs", "Hbegf", "Erbb2", "ECs", 3, 1, 4.540576, 13.621728, 0.1688497018860685, 0.1688497018860685, 2, 0.6666666666666666, 1.720171333333333, 5.160514, 0.1961456356393658, 0.1961456356393658, 7.810568672021333, 70.29511804819199, 0.03311913210396033, 0.03311913210396032),
    @("sCs", "Hbegf", "Erbb2", "FAPs", 3, 1, 4.540576, 13.621728, 0.1688497018860685, 0.1688497018860685, 3, 1, 4.076388666666666, 12.229166, 0.464817562438416, 0.464817562438416, 18.50915254653867, 166.582372918848, 0.07848430684913558, 0.07848430684913557),
    @("sCs", "Hbegf", "Erbb2", "sCs", 3, 1, 4.540576, 13.621728, 0.1688497018860685, 0.1688497018860685, 3, 1, 2.973308, 8.919924, 0.3390368019222182, 0.3390368019222182, 13.500530945408, 121.504778508672, 0.05724626293297261, 0.05724626293297259),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowData[$j]
    }
}